$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("Alt")
$src.Copy($null, $src)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "MySingleFunction"
foreach ($s in $wb.Worksheets) {
  Write-Output $s.Name
}
